$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the two samples that were removed from the dataset ---
# "RM 232" (row 26) and "SC 92" (originally row 28) are no longer present;
# deleting row 26 first shifts "SC 92" up to row 27, so deleting row 27
# next removes it too. Everything below cascades up by two rows.
$ws.Rows("26").Delete()
$ws.Rows("27").Delete()

# --- Re-roll which cells are "missing" for the remaining samples ---
# A numeric Value assignment fills in a previously-missing reading; the
# "'" + Style="Normal" pair writes a blank (empty-text) placeholder cell
# into a spot that used to hold a reading, without leaving a stray
# quote-prefix style behind.

# RM 2
$ws.Range("C2").Value = 14.9

# RM 8
$ws.Range("F3").Value = "'"
$ws.Range("F3").Style = "Normal"

# RM 9
$ws.Range("F4").Value = 17.97

# RM 21
$ws.Range("C6").Value = "'"
$ws.Range("C6").Style = "Normal"

# RM 38
$ws.Range("F8").Value = "'"
$ws.Range("F8").Style = "Normal"

# RM 42
$ws.Range("F9").Value = "'"
$ws.Range("F9").Style = "Normal"

# RM 81
$ws.Range("C12").Value = 12.5

# RM 90
$ws.Range("C14").Value = "'"
$ws.Range("C14").Style = "Normal"

# RM 95
$ws.Range("F15").Value = 16.2

# RM 120
$ws.Range("F18").Value = 18.35

# RM 125
$ws.Range("F19").Value = "'"
$ws.Range("F19").Style = "Normal"

# RM 134
$ws.Range("C20").Value = 12.5

# RM 135
$ws.Range("C21").Value = 12.7

# RM 138
$ws.Range("F22").Value = "'"
$ws.Range("F22").Style = "Normal"

# RM 140
$ws.Range("C23").Value = "'"
$ws.Range("C23").Style = "Normal"
$ws.Range("F23").Value = 16.48

# RM 142a
$ws.Range("C24").Value = "'"
$ws.Range("C24").Style = "Normal"

# RM 145
$ws.Range("F25").Value = 16.6

# SC 5 (now row 26)
$ws.Range("B26").Value = -20.2

# SC 101 (now row 27)
$ws.Range("B27").Value = "'"
$ws.Range("B27").Style = "Normal"
$ws.Range("F27").Value = "'"
$ws.Range("F27").Style = "Normal"

# SC 105 (now row 28)
$ws.Range("B28").Value = "'"
$ws.Range("B28").Style = "Normal"

# SC 119 (now row 29)
$ws.Range("B29").Value = -19.5

# SC 120 (now row 30)
$ws.Range("B30").Value = -19.7

# SC 132 (now row 31)
$ws.Range("B31").Value = "'"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = 15.3

# SC 193 (now row 32)
$ws.Range("B32").Value = "'"
$ws.Range("B32").Style = "Normal"

# SC 232 (now row 33)
$ws.Range("C33").Value = 10.4
